$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ Row=133; B=1994; E="Establishment"; F="Technical School at Chhauni, Kathmandu" },
    @{ Row=134; B=2014; E="Establishment"; F="Agriculture School in Minbhawan to train JT/JTAs" },
    @{ Row=135; B=2017; E="Establishment"; F="Jagadamba College of Agriculture and Research Institute at Shreemahal" },
    @{ Row=136; B=2025; E="Conversion"; F="School of Agriculture at Minbhawan transformed into Maha-vidyalaya and I.Sc. Agriculture program started. The program was launched in 2027 Shrawan 26." },
    @{ Row=137; B=2029; E="Conversion"; F="“Maha-vidyalaya” established in 2025 converted into IAAS with affiliation of Tribhuwan University and transferred to Jagadamba building at Pulchowk" },
    @{ Row=138; B=2030; E="Conversion"; F="IAAS relocated to Rampur, Chitwan" },
    @{ Row=139; B=2034; E="Establishment"; F="B.Sc Agriculture program started in Rampur, Chitwan" },
    @{ Row=140; B=2059; E="Establishment"; F="PhD program started in IAAS, TU" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

$ws.Range("D140").Select()

